$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 38 (weekly update adds a new observation) ---------------
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 11
$ws.Range("B38").Value = "Vega Monumental Concepción"
$ws.Range("C38").Value = "Bíobío"
$ws.Range("D38").Value = 44876
$ws.Range("E38").Value = 8
$ws.Range("F38").Value = 100112001
$ws.Range("G38").Value = "Berenjena"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 180
$ws.Range("K38").Value = 19000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = 19444
$ws.Range("N38").Value = "`$/caja 60 unidades"
$ws.Range("O38").Value = "Provincia de Huasco"
$ws.Range("P38").Value = 324
$ws.Range("Q38").Value = 60
$ws.Range("R38").Value = "Hortaliza"

# --- Insert new row 124 (second new weekly observation) ---------------------
$ws.Rows.Item(124).Insert()

$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Vega Monumental Concepción"
$ws.Range("C124").Value = "Bíobío"
$ws.Range("D124").Value = 44946
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100112001
$ws.Range("G124").Value = "Berenjena"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 180
$ws.Range("K124").Value = 19000
$ws.Range("L124").Value = 20000
$ws.Range("M124").Value = 19444
$ws.Range("N124").Value = "`$/caja 60 unidades"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 324
$ws.Range("Q124").Value = 60
$ws.Range("R124").Value = "Hortaliza"
